$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fix the misaligned data: the "titulo" (D) column values had slipped out
#    of sync with their "Grupo de Datos" (A/C) pairing by one row. Rotate the
#    titles down by one row (row 2's old title becomes row 8's, everything
#    else shifts down) and re-sync the A/C group columns with the title that
#    now lands on each row.
# ---------------------------------------------------------------------------

$ws.Range("A2").Value = "Grupo de Datos 2"
$ws.Range("C2").Value = "Dato grupo 2"
$ws.Range("D2").Value = "Equipo veterano da un gran espectaculo"

$ws.Range("D3").Value = "U.S. Robotics presenta hallazgo"

$ws.Range("D4").Value = "Se presenta el nuevo teléfono móvil en evento"

$ws.Range("D5").Value = "Se mejora la conducción autónoma de vehículos"

$ws.Range("A6").Value = "Grupo de Datos 1"
$ws.Range("C6").Value = "Dato grupo 1"
$ws.Range("D6").Value = "Fuccia OS sacude al mundo"

$ws.Range("D7").Value = "Tenemos campeona del mundial de volleiball"

# Row 8 keeps the same values ("Grupo de Datos 2" / "Dato grupo 2" /
# "Equipo veterano da un gran espectaculo") - only its formatting changes
# below.

# ---------------------------------------------------------------------------
# 2) Normalize formatting:
#    - The shared date numeric format is corrected from a date+time pattern
#      to a plain date pattern, and is now painted with the same yellow
#      highlight used by the rest of the data columns.
#    - fecha_publicacion (E) for the six "real" data rows (2-7) now uses that
#      highlighted date format.
#    - Row 8 (previously unstyled / inconsistently styled) is normalized to
#      match the look of the other data rows, with its date cell reverting to
#      the plain yellow fill (no special date format) used elsewhere.
# ---------------------------------------------------------------------------

$ws.Range("E2:E7").NumberFormat = "yyyy-mm-dd"
$ws.Range("E2:E7").Interior.Color = 65535

$ws.Range("A8").Interior.Color = 32768
$ws.Range("B8").Interior.Color = 8421504
$ws.Range("C8").Interior.Color = 65535
$ws.Range("D8").Interior.Color = 65535

# E8 should land on the exact same "plain yellow fill / General number
# format" look already used by C8/D8 - copy that format over instead of
# setting properties individually, so no redundant style gets created.
$ws.Range("D8").Copy()
$ws.Range("E8").PasteSpecial(-4122)
$excel.CutCopyMode = $false
